# Updated the add in to include the Log File information.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("config")

# Selection moved from C23 to C14
[void]$ws.Range("C14").Select()

# Simulator iterations changed
$ws.Range("B5").Value = 1209600

# Row 10: "template" -> "Template", and the stray "yes" (D10) is removed
$ws.Range("A10").Value = "Template"
$ws.Range("D10").ClearContents()

# Data validations: rebuild in the target order (A11, B9, A12, A10)
# A12's old validation is removed and re-added later (with an updated list)
# so that the declaration order matches the authored workbook.
$ws.Range("A12").Validation.Delete()

$ws.Range("B9").Validation.Add(3, 1, 1, '"Yes, No"')
$ws.Range("B9").Validation.ShowError = $false

$ws.Range("A12").Validation.Add(3, 1, 1, '",,,,,Community Name,Template,,input,output"')
$ws.Range("A12").Validation.ShowError = $false

$ws.Range("A10").Validation.Add(3, 1, 1, '",,,,,Community Name,,,input,output,Template"')
$ws.Range("A10").Validation.ShowError = $false
